$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header fields for consistency
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update selection to reflect the edited cells
$ws.Range("K1:L1").Select()
